$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while guaranteeing it stays text, even when the
# string looks like a number (Excel would otherwise parse it numerically).
# We flip the cell to a text number-format just long enough to assign the
# value, then restore the cells original Style object so formatting is
# unchanged.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "67.558.71"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "3.869.59"
$ws.Range("E3").Value = "  +0.61%  "
Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue $ws.Range("D5") "461.81"
$ws.Range("E5").Value = "  +9.22%  "
Set-TextValue $ws.Range("D6") "147.15"
$ws.Range("E6").Value = "  +14.18%  "
Set-TextValue $ws.Range("D7") "0.626"
$ws.Range("E7").Value = "  +2.98%  "
$ws.Range("E9").Value = "  +3.84%  "
Set-TextValue $ws.Range("D10") "0.156"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("E11").Value = "  -5.77%  "
Set-TextValue $ws.Range("D12") "43.97"
$ws.Range("E12").Value = "  +7.93%  "
Set-TextValue $ws.Range("D13") "10.42"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Value = "4.463.86"
$ws.Range("E14").Value = "  +0.31%  "
Set-TextValue $ws.Range("D15") "14.83"
$ws.Range("E15").Value = "  -4.84%  "
$ws.Range("D16").Value = "3.873.43"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("E17").Value = "  -0.17%  "
Set-TextValue $ws.Range("D18") "20.08"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("E19").Value = "  +7.17%  "
$ws.Range("D20").Value = "67.547.73"
$ws.Range("E20").Value = "  +0.65%  "
Set-TextValue $ws.Range("D21") "427.66"
$ws.Range("E21").Value = "  +4.45%  "
Set-TextValue $ws.Range("D22") "14.90"
$ws.Range("E22").Value = "  -0.11%  "
Set-TextValue $ws.Range("D23") "3.28"
$ws.Range("E23").Value = "  +8.28%  "
Set-TextValue $ws.Range("D24") "87.02"
$ws.Range("E24").Value = "  +3.41%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D25") "3.55"
$ws.Range("E25").Value = "  +9.46%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D26") "10.53"
$ws.Range("E26").Value = "  +11.67%  "
Set-TextValue $ws.Range("D27") "37.65"
$ws.Range("E27").Value = "  +0.32%  "
Set-TextValue $ws.Range("D28") "10.06"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  +2.31%  "
Set-TextValue $ws.Range("D30") "756.99"
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("E31").Value = "  +12.03%  "
$ws.Range("E32").Value = "  +5.41%  "
Set-TextValue $ws.Range("D33") "2.75"
$ws.Range("E33").Value = "  -0.97%  "
Set-TextValue $ws.Range("D34") "43.64"
$ws.Range("E34").Value = "  +13.72%  "
$ws.Range("E35").Value = "  +7.40%  "
Set-TextValue $ws.Range("D36") "57.57"
$ws.Range("E36").Value = "  +3.81%  "
Set-TextValue $ws.Range("D37") "5.56"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("E38").Value = "  +0.01%  "
Set-TextValue $ws.Range("D39") "0.0479"
$ws.Range("E39").Value = "  +4.80%  "
Set-TextValue $ws.Range("D40") "0.357"
$ws.Range("E40").Value = "  +12.36%  "
$ws.Range("E41").Value = "  +2.25%  "
Set-TextValue $ws.Range("D42") "2.67"
$ws.Range("E42").Value = "  +16.48%  "
$ws.Range("D43").Value = "0.0₃0678"
$ws.Range("E43").Value = "  -6.51%  "
$ws.Range("E44").Value = "  +5.19%  "
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("E48").Value = "  +7.87%  "
$ws.Range("E49").Value = "  +4.56%  "
Set-TextValue $ws.Range("D50") "144.61"
$ws.Range("E50").Value = "  +2.93%  "
Set-TextValue $ws.Range("D51") "2.89"
$ws.Range("E51").Value = "  +2.94%  "
